# Fruta / hortaliza, semanal
# Permute the per-row observation data (Fecha, Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion, Origen, Precio $/Kg, Kg/unidad) across rows 2-10.
# Columns A, B, C, E, F, G, H, I, J, K, L are identical for every row and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D, M, N, O, P, Q, R, S, T
$data = @{
    2  = @{ D = 44594; M = 120; N = 2500; O = 2800; P = 2650; Q = "$/bandeja 2 kilos"; R = "Provincia de Linares";  S = 1325; T = 2 }
    3  = @{ D = 44539; M = 200; N = 3800; O = 4000; P = 3900; Q = "$/bandeja 2 kilos"; R = "Región del Maule";      S = 1950; T = 2 }
    4  = @{ D = 44540; M = 240; N = 3500; O = 3800; P = 3650; Q = "$/bandeja 2 kilos"; R = "Región del Maule";      S = 1825; T = 2 }
    5  = @{ D = 44596; M = 120; N = 2500; O = 2700; P = 2600; Q = "$/bandeja 2 kilos"; R = "Provincia de Linares";  S = 1300; T = 2 }
    6  = @{ D = 44187; M = 80;  N = 2800; O = 3000; P = 2900; Q = "$/bandeja 2 kilos"; R = "Provincia de Linares";  S = 1450; T = 2 }
    7  = @{ D = 44187; M = 65;  N = 1400; O = 1500; P = 1446; Q = "$/envase 1 kilo";   R = "Provincia de Diguillín"; S = 1446; T = 1 }
    8  = @{ D = 44174; M = 150; N = 3700; O = 3800; P = 3747; Q = "$/bandeja 2 kilos"; R = "Provincia de Linares";  S = 1874; T = 2 }
    9  = @{ D = 44181; M = 65;  N = 3600; O = 3800; P = 3692; Q = "$/bandeja 2 kilos"; R = "Provincia de Diguillín"; S = 1846; T = 2 }
    10 = @{ D = 44181; M = 80;  N = 1800; O = 2000; P = 1875; Q = "$/envase 1 kilo";   R = "Provincia de Diguillín"; S = 1875; T = 1 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
    $ws.Range("T$row").Value = $vals.T
}
